$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" '56.648.88'
Set-TextValue "E2" '  -2.61%  '
Set-TextValue "D3" '2.987.22'
Set-TextValue "E4" '  +0.05%  '
Set-TextValue "D5" '497.53'
Set-TextValue "E5" '  -5.11%  '
Set-TextValue "D6" '134.42'
Set-TextValue "E6" '  -0.35%  '
Set-TextValue "E7" '  -0.01%  '
Set-TextValue "D8" '2.985.50'
Set-TextValue "E8" '  -4.77%  '
Set-TextValue "E9" '  -4.24%  '
Set-TextValue "D10" '7.27'
Set-TextValue "E10" '  +0.10%  '
Set-TextValue "E11" '  -3.81%  '
Set-TextValue "E12" '  -7.73%  '
Set-TextValue "E13" '  -0.49%  '
Set-TextValue "D14" '3.493.77'
Set-TextValue "E14" '  -4.91%  '
Set-TextValue "D15" '24.78'
Set-TextValue "E15" '  -3.22%  '
Set-TextValue "D16" '56.626.87'
Set-TextValue "E16" '  -2.65%  '
Set-TextValue "E17" '  -3.41%  '
Set-TextValue "D18" '2.982.95'
Set-TextValue "E18" '  -4.98%  '
Set-TextValue "E19" '  +0.31%  '
Set-TextValue "D20" '12.31'
Set-TextValue "E20" '  -5.82%  '
Set-TextValue "E21" '  -2.38%  '
Set-TextValue "D22" '326.00'
Set-TextValue "E22" '  -5.27%  '
Set-TextValue "D23" '1.00'
Set-TextValue "E23" '  +0.12%  '
Set-TextValue "E24" '  -8.43%  '
Set-TextValue "D25" '61.29'
Set-TextValue "E25" '  -10.55%  '
Set-TextValue "D26" '0.996'
Set-TextValue "E26" '  -0.36%  '
Set-TextValue "E27" '  -3.77%  '
Set-TextValue "D28" '0.0₃0911'
Set-TextValue "E28" '  -4.83%  '
Set-TextValue "D29" '0.998'
Set-TextValue "E29" '  -0.03%  '
Set-TextValue "D30" '6.49'
Set-TextValue "E30" '  -4.62%  '
Set-TextValue "D31" '6.81'
Set-TextValue "E31" '  -0.70%  '
Set-TextValue "D32" '1.17'
Set-TextValue "E32" '  -4.02%  '
Set-TextValue "E33" '  -6.70%  '
Set-TextValue "E34" '  -7.05%  '
Set-TextValue "D35" '154.25'
Set-TextValue "E35" '  -1.73%  '
Set-TextValue "D36" '4.47'
Set-TextValue "E36" '  -6.76%  '
Set-TextValue "E37" '  -6.85%  '
Set-TextValue "D38" '5.60'
Set-TextValue "E38" '  -9.99%  '
Set-TextValue "D39" '0.0675'
Set-TextValue "E39" '  -2.23%  '
Set-TextValue "D40" '23.47'
Set-TextValue "E40" '  -3.76%  '
Set-TextValue "D41" '3.018.40'
Set-TextValue "E41" '  -4.71%  '
Set-TextValue "D42" '36.68'
Set-TextValue "E42" '  -9.19%  '
Set-TextValue "E43" '  +0.06%  '
Set-TextValue "E44" '  -6.77%  '
Set-TextValue "D45" '0.637'
Set-TextValue "E45" '  -7.83%  '
Set-TextValue "E46" '  -2.18%  '
Set-TextValue "D47" '2.200.22'
Set-TextValue "E47" '  -2.59%  '
Set-TextValue "E48" '  -8.79%  '
Set-TextValue "E49" '  +6.41%  '
Set-TextValue "E50" '  +1.73%  '
Set-TextValue "D51" '5.71'
Set-TextValue "E51" '  -7.68%  '
